$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.610.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '''3.522.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '''597.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = '''143.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.05%  '
$ws.Range("D7").Value = '''3.520.44'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '''0.503'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("D12").Value = '''0.404'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").Value = '''4.122.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("E14").Value = '  -4.05%  '
$ws.Range("D15").Value = '''28.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.33%  '
$ws.Range("D16").Value = '''3.501.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("E17").Value = '  +1.33%  '
$ws.Range("D18").Value = '''65.635.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").Value = '''10.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.49%  '
$ws.Range("D20").Value = '''6.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.33%  '
$ws.Range("D21").Value = '''14.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("D22").Value = '''415.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").Value = '''0.599'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("D24").Value = '''77.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.60%  '
$ws.Range("D25").Value = '''3.663.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '''0.0000115'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.01%  '
$ws.Range("D28").Value = '''2.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.92%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''8.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.47%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '''7.73'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '''3.521.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").Value = '''0.152'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").Value = '''24.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.90%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").Value = '''7.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.49%  '
$ws.Range("D37").Value = '''1.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.57%  '
$ws.Range("D38").Value = '''174.44'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").Value = '''5.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.70%  '
$ws.Range("D40").Value = '''1.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.81%  '
$ws.Range("D41").Value = '''0.0820'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.42%  '
$ws.Range("D42").Value = '''5.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").Value = '''0.856'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.61%  '
$ws.Range("D44").Value = '''45.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.93%  '
$ws.Range("E45").Value = '  -8.18%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("D47").Value = '''2.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.41%  '
$ws.Range("D48").Value = '''7.09'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").Value = '''22.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.24%  '
$ws.Range("E50").Value = '  -8.57%  '
$ws.Range("D51").Value = '''22.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -9.22%  '
